# Realestate Update resale numbers 2023-05-31 21:38
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start the new row (13) as a copy of the previous row (12) so that the
# text-like columns (Date/Time/Weekday/Week) keep their original text
# representation instead of being reinterpreted (e.g. dates auto-converted
# to serial numbers).
$ws.Range("A12:T12").Copy()
$ws.Range("A13").PasteSpecial(-4104)

# Now overwrite the cells that actually changed for the new entry.
$ws.Cells.Item(13, 2).Value = "21:36:21"   # Time

$ws.Cells.Item(13, 5).Value = 120271    # Beijing
$ws.Cells.Item(13, 6).Value = 133519    # Guangzhou
$ws.Cells.Item(13, 7).Value = 158974    # Suzhou
$ws.Cells.Item(13, 8).Value = 130957    # Hangzhou
$ws.Cells.Item(13, 9).Value = 174859    # Nanjing
$ws.Cells.Item(13, 10).Value = 113878   # Xi_an
$ws.Cells.Item(13, 11).Value = 198929   # Chengdu
$ws.Cells.Item(13, 12).Value = 220338   # Chongqing
$ws.Cells.Item(13, 13).Value = 172243   # Tianjin
$ws.Cells.Item(13, 14).Value = 120038   # Hefei
$ws.Cells.Item(13, 15).Value = 38729    # Fuzhou
$ws.Cells.Item(13, 16).Value = 34925    # Xiamen
$ws.Cells.Item(13, 17).Value = 50654    # Changsha
$ws.Cells.Item(13, 19).Value = 36961    # Shenzhen
